$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 222, shifting existing rows (222-236) down to (223-237).
$ws.Range("A222:T222").Insert("xlShiftDown")

# Populate the newly inserted row 222 with this week's data point.
$ws.Range("A222").Value = 10
$ws.Range("B222").Value = "Vega Modelo de Temuco"
$ws.Range("C222").Value = "La Araucanía"
$ws.Range("D222").Value = 45223
$ws.Range("E222").Value = 9
$ws.Range("F222").Value = "Fruta"
$ws.Range("G222").Value = 100107
$ws.Range("H222").Value = "Otros"
$ws.Range("I222").Value = 100107002
$ws.Range("J222").Value = "Chirimoya"
$ws.Range("K222").Value = "Cultivar IV Región"
$ws.Range("L222").Value = "Primera"
$ws.Range("M222").Value = 55
$ws.Range("N222").Value = 2600
$ws.Range("O222").Value = 2600
$ws.Range("P222").Value = 2600
$ws.Range("Q222").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R222").Value = "Provincia de Limarí"
$ws.Range("S222").Value = 2600
$ws.Range("T222").Value = 1
